$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# "create update request methods"
# ---------------------------------------------------------------------

# 1) Add "delete video comments" / 35 to the Delete table (G/H columns),
#    continuing the existing list (G1:H5 -> G6:H6), matching the style
#    used by the rest of that table.
$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G6:H6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("G6").Value = "delete video comments"
$ws.Range("H6").Value = 35

# 2) Add three new rows to the Update table (J/K columns), continuing the
#    existing list (J1:K10 -> J11:K13), matching the style used by the
#    rest of that table.
$ws.Range("J10:K10").Copy() | Out-Null
$ws.Range("J11:K13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("J11").Value = "update video views"
$ws.Range("K11").Value = 410

$ws.Range("J12").Value = "update video likes"
$ws.Range("K12").Value = 411

$ws.Range("J13").Value = "update video deslike "
$ws.Range("K13").Value = 412

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Cosmetic follow-up: selection moved to K14 and columns were resized
# (best effort - the headless column-width model here quantizes widths
# to a coarser grid than desktop Excel, so these land on the nearest
# reachable width rather than the exact fractional value).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.746666666666666
$ws.Columns.Item(2).ColumnWidth = 3.08
$ws.Columns.Item(3).ColumnWidth = 1.08
$ws.Columns.Item(6).ColumnWidth = 1.413333333333333
$ws.Columns.Item(7).ColumnWidth = 18.746666666666666
$ws.Columns.Item(8).ColumnWidth = 2.08
$ws.Columns.Item(9).ColumnWidth = 1.58
$ws.Columns.Item(12).ColumnWidth = 1.58

$ws.Range("K14").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 3
} catch {}

Write-Host "done"
